# Updates the "cryptos" price table on Sheet1 (rows 2-51) so that the
# Price (column D) and Volume(1h) (column E) columns reflect the latest
# scrape, and re-ranks a handful of coins whose relative order changed
# (rows 20/21, 29/30, 32/33, 49/50/51 swap identities).
#
# $CellUpdates is an ordered list of (cell, new text value) pairs taken
# directly from the target diff; it is applied in document order so that
# the B/C/D/E cells for a re-ranked row are all rewritten together.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$CellUpdates = @(
    @{ Cell = 'D2'; Value = '42.198.93' }
    @{ Cell = 'E2'; Value = '  -1.91%  ' }
    @{ Cell = 'D3'; Value = '2.249.36' }
    @{ Cell = 'E3'; Value = '  -1.97%  ' }
    @{ Cell = 'E4'; Value = '  +0.04%  ' }
    @{ Cell = 'D5'; Value = '247.14' }
    @{ Cell = 'E5'; Value = '  -1.95%  ' }
    @{ Cell = 'D6'; Value = '0.634' }
    @{ Cell = 'E6'; Value = '  -1.28%  ' }
    @{ Cell = 'D7'; Value = '77.20' }
    @{ Cell = 'E7'; Value = '  +3.83%  ' }
    @{ Cell = 'E8'; Value = '  +0.01%  ' }
    @{ Cell = 'D9'; Value = '0.625' }
    @{ Cell = 'E9'; Value = '  -3.51%  ' }
    @{ Cell = 'D10'; Value = '41.67' }
    @{ Cell = 'E10'; Value = '  +5.92%  ' }
    @{ Cell = 'E11'; Value = '  -2.90%  ' }
    @{ Cell = 'E12'; Value = '  -4.57%  ' }
    @{ Cell = 'E13'; Value = '  -3.19%  ' }
    @{ Cell = 'D14'; Value = '2.585.50' }
    @{ Cell = 'E14'; Value = '  -1.96%  ' }
    @{ Cell = 'D15'; Value = '14.85' }
    @{ Cell = 'E15'; Value = '  -3.52%  ' }
    @{ Cell = 'D16'; Value = '0.862' }
    @{ Cell = 'E16'; Value = '  -1.58%  ' }
    @{ Cell = 'D17'; Value = '2.231.14' }
    @{ Cell = 'E17'; Value = '  -2.59%  ' }
    @{ Cell = 'D18'; Value = '42.101.38' }
    @{ Cell = 'E18'; Value = '  -1.96%  ' }
    @{ Cell = 'D19'; Value = '0.0₃0988' }
    @{ Cell = 'E19'; Value = '  -2.72%  ' }
    @{ Cell = 'B20'; Value = 'Litecoin' }
    @{ Cell = 'C20'; Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc' }
    @{ Cell = 'D20'; Value = '72.02' }
    @{ Cell = 'E20'; Value = '  -0.96%  ' }
    @{ Cell = 'B21'; Value = 'Uniswap' }
    @{ Cell = 'C21'; Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni' }
    @{ Cell = 'D21'; Value = '6.12' }
    @{ Cell = 'E21'; Value = '  -2.90%  ' }
    @{ Cell = 'E22'; Value = '  +1.22%  ' }
    @{ Cell = 'D23'; Value = '232.13' }
    @{ Cell = 'E23'; Value = '  -2.44%  ' }
    @{ Cell = 'E24'; Value = '  -0.02%  ' }
    @{ Cell = 'D25'; Value = '11.36' }
    @{ Cell = 'E25'; Value = '  -2.53%  ' }
    @{ Cell = 'E26'; Value = '  -7.68%  ' }
    @{ Cell = 'E27'; Value = '  -4.92%  ' }
    @{ Cell = 'D28'; Value = '7.47' }
    @{ Cell = 'E28'; Value = '  +17.32%  ' }
    @{ Cell = 'B29'; Value = 'Monero' }
    @{ Cell = 'C29'; Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr' }
    @{ Cell = 'D29'; Value = '169.91' }
    @{ Cell = 'E29'; Value = '  +1.60%  ' }
    @{ Cell = 'B30'; Value = 'Toncoin' }
    @{ Cell = 'C30'; Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton' }
    @{ Cell = 'D30'; Value = '2.10' }
    @{ Cell = 'E30'; Value = '  -2.02%  ' }
    @{ Cell = 'D31'; Value = '20.83' }
    @{ Cell = 'E31'; Value = '  -1.36%  ' }
    @{ Cell = 'B32'; Value = 'Hedera' }
    @{ Cell = 'C32'; Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar' }
    @{ Cell = 'D32'; Value = '0.0838' }
    @{ Cell = 'E32'; Value = '  -0.02%  ' }
    @{ Cell = 'B33'; Value = 'InjectiveProtocol' }
    @{ Cell = 'C33'; Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj' }
    @{ Cell = 'D33'; Value = '32.95' }
    @{ Cell = 'E33'; Value = '  +5.74%  ' }
    @{ Cell = 'E34'; Value = '  -5.15%  ' }
    @{ Cell = 'E35'; Value = '  -1.38%  ' }
    @{ Cell = 'D36'; Value = '4.53' }
    @{ Cell = 'E36'; Value = '  -1.70%  ' }
    @{ Cell = 'D37'; Value = '4.94' }
    @{ Cell = 'E37'; Value = '  +2.60%  ' }
    @{ Cell = 'D39'; Value = '14.26' }
    @{ Cell = 'E39'; Value = '  +1.51%  ' }
    @{ Cell = 'D40'; Value = '5.91' }
    @{ Cell = 'E40'; Value = '  +0.06%  ' }
    @{ Cell = 'E41'; Value = '  -6.76%  ' }
    @{ Cell = 'D42'; Value = '113.27' }
    @{ Cell = 'E42'; Value = '  +8.77%  ' }
    @{ Cell = 'E43'; Value = '  -6.85%  ' }
    @{ Cell = 'D44'; Value = '61.12' }
    @{ Cell = 'E44'; Value = '  -1.87%  ' }
    @{ Cell = 'D45'; Value = '8.68' }
    @{ Cell = 'E45'; Value = '  -5.07%  ' }
    @{ Cell = 'D46'; Value = '0.0995' }
    @{ Cell = 'E46'; Value = '  -3.75%  ' }
    @{ Cell = 'E47'; Value = '  -0.63%  ' }
    @{ Cell = 'E48'; Value = '  -3.52%  ' }
    @{ Cell = 'B49'; Value = 'TrustWalletToken' }
    @{ Cell = 'C49'; Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt' }
    @{ Cell = 'D49'; Value = '1.17' }
    @{ Cell = 'E49'; Value = '  -1.09%  ' }
    @{ Cell = 'B50'; Value = 'FTXToken' }
    @{ Cell = 'C50'; Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt' }
    @{ Cell = 'D50'; Value = '4.33' }
    @{ Cell = 'E50'; Value = '  -12.02%  ' }
    @{ Cell = 'B51'; Value = 'WOONetwork' }
    @{ Cell = 'C51'; Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo' }
    @{ Cell = 'D51'; Value = '0.446' }
    @{ Cell = 'E51'; Value = '  +15.90%  ' }

)

# Many of the new Price values are plain decimal-looking strings
# (e.g. "247.14", "0.634"). Excel auto-coerces a plain numeric string
# assigned to Range.Value into a real number, but every Price/Volume cell
# on this sheet is stored as text. Force text ("@") number format before
# writing any such value, then restore the cell's original (default)
# style afterwards so no visible formatting change is left behind - only
# the underlying value stays text, matching the source workbook.
function Looks-Numeric($s) {
    return ($s -match '^[+-]?\d+(\.\d+)?$')
}

foreach ($u in $CellUpdates) {
    if (Looks-Numeric $u.Value) {
        $ws.Range($u.Cell).NumberFormat = "@"
    }
}

foreach ($u in $CellUpdates) {
    $ws.Range($u.Cell).Value = $u.Value
}

foreach ($u in $CellUpdates) {
    if (Looks-Numeric $u.Value) {
        $ws.Range($u.Cell).Style = "Normal"
    }
}
